# Loan RBI, Variable Instalments
# Insert a new column on the "Repayment Schedule" sheet (before column N) to make
# room for an extra "Variable Instalment"-style column, then switch the active
# sheet/selection from "Transactions" to "Repayment Schedule".

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before N; this shifts N->O, O->P, P->Q, etc. and keeps
# everything else (values/styles) intact, matching the original columns moving
# one slot to the right.
$wsSchedule.Columns("N").Insert()

# The new column N should carry the same width as column M (its left neighbour).
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab, with R7 selected, and drop the
# previous selection/active-tab state that lived on "Transactions".
$wsSchedule.Activate()
$wsSchedule.Range("R7").Select()
